# Update FHIR ValueSet metadata from ibm.com/Alvearie to linuxforhealth.org/LinuxForHealth
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$metaSheet = $wb.Worksheets.Item("Metadata")

# URL (B2)
$metaSheet.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/hipaa-relationship"

# Version (B3)
$metaSheet.Range("B3").Value = "8.0.0"

# Date (B8)
$metaSheet.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher (B9)
$metaSheet.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet 2: Include from HIPAA Relationsh ---
$codeSheet = $wb.Worksheets.Item("Include from HIPAA Relationsh")

# System URI (B4)
$codeSheet.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/hipaa-relationship"
